$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update runtime (column F) values for the TELE database rows for every
# detector, reflecting the new MATLAB runtime measurements.
$ws.Range("F30").Value = 473.2589721679688
$ws.Range("F31").Value = 507.0305585861206
$ws.Range("F32").Value = 672.3466515541077
$ws.Range("F33").Value = 727.0430326461792
$ws.Range("F34").Value = 746.0289001464844
$ws.Range("F35").Value = 864.7780656814575
$ws.Range("F36").Value = 880.8329105377197
$ws.Range("F68").Value = 101.2722253799438
$ws.Range("F69").Value = 325.6262540817261
$ws.Range("F70").Value = 408.578073978424
$ws.Range("F71").Value = 432.3427677154541
$ws.Range("F100").Value = 9.634971618652344
$ws.Range("F101").Value = 11.90927028656006
$ws.Range("F102").Value = 20.73484659194946
$ws.Range("F103").Value = 42.24205017089844
$ws.Range("F104").Value = 74.79268312454224
$ws.Range("F105").Value = 101.5661597251892
$ws.Range("F106").Value = 112.6465797424316
$ws.Range("F135").Value = 333.3365917205811
$ws.Range("F136").Value = 335.7885599136353
$ws.Range("F137").Value = 352.6906967163086
$ws.Range("F138").Value = 449.8480558395386
$ws.Range("F139").Value = 468.2725667953491
$ws.Range("F140").Value = 545.2590465545653
$ws.Range("F141").Value = 584.355354309082
$ws.Range("F170").Value = 33.15544128417969
$ws.Range("F171").Value = 36.6207480430603
$ws.Range("F172").Value = 50.86004734039307
$ws.Range("F173").Value = 56.82849884033203
$ws.Range("F174").Value = 65.67549705505371
$ws.Range("F175").Value = 68.92014741897583
$ws.Range("F176").Value = 70.05214691162109
$ws.Range("F205").Value = 16.72506332397461
$ws.Range("F206").Value = 18.6953067779541
$ws.Range("F207").Value = 28.81044149398804
$ws.Range("F208").Value = 37.89198398590088
$ws.Range("F209").Value = 52.61421203613281
$ws.Range("F210").Value = 222.1382975578305
$ws.Range("F211").Value = 349.463939666748
$ws.Range("F241").Value = 3.598880767822266
$ws.Range("F242").Value = 8.295178413391113
$ws.Range("F243").Value = 9.53376293182373
$ws.Range("F244").Value = 10.05858182907104
$ws.Range("F245").Value = 15.08040428161621
$ws.Range("F246").Value = 15.7923698425293
$ws.Range("F275").Value = 4.53639030456543
$ws.Range("F276").Value = 6.079626083374023
$ws.Range("F277").Value = 8.040904998779297
$ws.Range("F278").Value = 9.339213371276855
$ws.Range("F279").Value = 14.96833562850952
$ws.Range("F280").Value = 18.06657314300537
$ws.Range("F281").Value = 19.02937889099121
$ws.Range("F310").Value = 40.66205024719238
$ws.Range("F311").Value = 43.96588802337646
$ws.Range("F312").Value = 59.2692494392395
$ws.Range("F313").Value = 74.21457767486572
$ws.Range("F314").Value = 88.94813060760498
$ws.Range("F315").Value = 146.2499022483825
$ws.Range("F316").Value = 170.6409454345703
$ws.Range("F345").Value = 5.041599273681641
$ws.Range("F346").Value = 5.772662162780762
$ws.Range("F347").Value = 7.713854312896729
$ws.Range("F348").Value = 7.993698120117188
$ws.Range("F349").Value = 9.370386600494385
$ws.Range("F350").Value = 16.01730585098267
$ws.Range("F351").Value = 16.48926734924316
$ws.Range("F382").Value = 3.662526607513428
$ws.Range("F383").Value = 8.046388626098633
$ws.Range("F384").Value = 8.676350116729736
$ws.Range("F385").Value = 10.00897884368896
$ws.Range("F386").Value = 10.29586791992188
$ws.Range("F416").Value = 12.79209852218628
$ws.Range("F417").Value = 31.01474046707153
$ws.Range("F418").Value = 39.22808170318604
$ws.Range("F419").Value = 72.36778736114502
$ws.Range("F420").Value = 191.738224029541
$ws.Range("F421").Value = 219.0546989440918
$ws.Range("F450").Value = 22.91059494018555
$ws.Range("F451").Value = 22.94224500656128
$ws.Range("F452").Value = 27.09120512008667
$ws.Range("F453").Value = 36.63945198059082
$ws.Range("F454").Value = 50.99308490753174
$ws.Range("F455").Value = 190.881597995758
$ws.Range("F456").Value = 207.3073387145996
$ws.Range("F485").Value = 12.30263710021973
$ws.Range("F486").Value = 13.17306756973267
$ws.Range("F487").Value = 21.68565988540649
$ws.Range("F488").Value = 23.93829822540283
$ws.Range("F489").Value = 25.16293525695801
$ws.Range("F490").Value = 30.24462461471557
$ws.Range("F491").Value = 32.64284133911133
$ws.Range("F520").Value = 102.527379989624
$ws.Range("F521").Value = 107.3583602905273
$ws.Range("F522").Value = 119.4368004798889
$ws.Range("F523").Value = 137.0326280593872
$ws.Range("F524").Value = 177.2399544715881
$ws.Range("F525").Value = 250.5567669868469
$ws.Range("F526").Value = 253.9701461791992
$ws.Range("F555").Value = 653.5487174987793
$ws.Range("F556").Value = 680.7798981666565
$ws.Range("F557").Value = 776.4614224433899
$ws.Range("F558").Value = 918.9643859863281
$ws.Range("F559").Value = 1015.517115592957
$ws.Range("F560").Value = 1063.24657201767
$ws.Range("F561").Value = 1097.378969192505
